$d = $word.ActiveDocument

# Locate the reference paragraph that ends with the "( <hyperlink> )" text
# (unique substring taken from the tracked hyperlink URL) and expand the
# found range to cover the whole paragraph so we can collapse to its end.
$anchor = $d.Content
$found = $anchor.Find.Execute("source=images", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph for Persona insertion"
}
$null = $anchor.Expand(4)
$anchor.Collapse(0)

# Build an insertion point Range positioned exactly at the start of the
# first trailing empty paragraph, then inject the new OOXML as siblings
# in front of it (InsertXML creates new paragraphs rather than disturbing
# the two empty paragraphs that must remain at the very end of the body).
$insertionPoint = $d.Range($anchor.Start, $anchor.Start)

$personaXml = @'
<w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Persona</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Career</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Siobhan, carer of a 65-year-old patient</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Siobhan is a 25-year-old carer who just finished college and is now doing her first job in a carer home. She been tasked to take care of John who is a 65-year-old suffering from dementia. It takes her almost 2 hours to commute from Kildare to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ardee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>in order to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> be able to work. As she </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>has to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> take good care of him, she (as the carer) sign john (as the patient) to the care tracker app as it easier for her to keep track of john.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">By signing him up he </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is able to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> let his carer know about his where abouts and possible problems when she isn't available/ around</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Users Details</w:t></w:r></w:p><w:p><w:r><w:t>Name: Siobhan Kerr</w:t></w:r></w:p><w:p><w:r><w:t>DOB: 10th September 1995</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Occupation: Carer  </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Carer Home: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ardee</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Care Home</w:t></w:r></w:p><w:p><w:r><w:t>Patient: John Stacy</w:t></w:r></w:p><w:p><w:r><w:t>Age: 55 years</w:t></w:r></w:p><w:p><w:r><w:t>Sickness: Dementia</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Siobhan's Goals</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>It ensures John that suffering with dementia is safe at all cost</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">Her goal with the care tracker app is to be able to monitor John knowing she </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>won’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> be around him most of the time so it would help her keep informed</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>By selecting the Care tracker App, John will be able to let Siobhan know his whereabouts without her having to worry about it.</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Siobhan is also able to look at John’s planner and see what he has planned for the day or the week.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Interface Requirements</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Allows to tracker</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Provide a planner</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Minimise obsolete screens making user interface nice and easy for both users</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Allow Carer to be able to accesses information</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Persona</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Client view</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Maire, newly retired labourer </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Marie is a 67-year-old newly retired labourer who stays in Cork after working in a factory for 40 years has been suffering from diabetes. Her Carer which is Philip is usually her stay at home carer but he decided to go and visit his family for the weekend and fears that Maire may be in danger as her vision and hearing has started to deuterate which impedes her daily. As Maire is aware, she hears about this app that her friend Mildred uses to keep in contact with her career. She hears </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>it’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> called Care Tracker App and it’s as simple to use with a 4-part navbar so if she needs to contact Phillip it’s at a touch of a button. She signs up and tells Phillip to set up one as the carer. This way she can keep in contact with Phillip and he can give her instruction on prescriptions </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Users Details</w:t></w:r></w:p><w:p><w:r><w:t>Name: Maire Stacy</w:t></w:r></w:p><w:p><w:r><w:t>DOB: 14 Feb 1951</w:t></w:r></w:p><w:p><w:r><w:t>Age:67 years</w:t></w:r></w:p><w:p><w:r><w:t>Address: 22 Manor Road</w:t></w:r></w:p><w:p><w:r><w:t>Sickness: Diabetes</w:t></w:r></w:p><w:p><w:r><w:t>Career: Phillip Mc Donald</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Maire Goals</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• It ensures that she can still receive instruction from Philip </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>in order to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> be able to know which prescription to take </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">• Her goal with the care tracker app is to be able to keep in contact with Phillip whilst </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>he’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> visiting home</w:t></w:r></w:p><w:p><w:r><w:t>• By selecting the Care tracker App, Maire will be able to let Phillip know her whereabouts without him having to worry about it.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Interface Requirements</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Allows to tracker</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Provide a planner</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Minimise obsolete screens making user interface nice and easy for both users</w:t></w:r></w:p><w:p><w:r><w:t>•</w:t></w:r><w:r><w:tab/><w:t>Allow Carer to be able to accesses information</w:t></w:r></w:p><w:p/>
'@

$null = $insertionPoint.InsertXML($personaXml)
